$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full rewrite of the sheet so spans/dimension metadata are recomputed cleanly
$ws.Cells.Clear()

# Re-apply header style (bold, centered, thin border) to header row and label column
$headerStyle = $ws.Range("B1:K1,A2:A24")
$headerStyle.Font.Bold = $true
$headerStyle.HorizontalAlignment = -4108
$headerStyle.VerticalAlignment = -4160
$headerStyle.Borders.LineStyle = 1

# Row 1 header labels (Q0..Q9)
$ws.Range("B1").Value2 = "Q0"
$ws.Range("C1").Value2 = "Q1"
$ws.Range("D1").Value2 = "Q2"
$ws.Range("E1").Value2 = "Q3"
$ws.Range("F1").Value2 = "Q4"
$ws.Range("G1").Value2 = "Q5"
$ws.Range("H1").Value2 = "Q6"
$ws.Range("I1").Value2 = "Q7"
$ws.Range("J1").Value2 = "Q8"
$ws.Range("K1").Value2 = "Q9"

# Column A row labels (quarterly dates)
$ws.Range("A2").Value2 = "2020-01-01 00:00:00_diff"
$ws.Range("A3").Value2 = "2020-04-01 00:00:00_diff"
$ws.Range("A4").Value2 = "2020-07-01 00:00:00_diff"
$ws.Range("A5").Value2 = "2020-10-01 00:00:00_diff"
$ws.Range("A6").Value2 = "2021-01-01 00:00:00_diff"
$ws.Range("A7").Value2 = "2021-04-01 00:00:00_diff"
$ws.Range("A8").Value2 = "2021-07-01 00:00:00_diff"
$ws.Range("A9").Value2 = "2021-10-01 00:00:00_diff"
$ws.Range("A10").Value2 = "2022-01-01 00:00:00_diff"
$ws.Range("A11").Value2 = "2022-04-01 00:00:00_diff"
$ws.Range("A12").Value2 = "2022-07-01 00:00:00_diff"
$ws.Range("A13").Value2 = "2022-10-01 00:00:00_diff"
$ws.Range("A14").Value2 = "2023-01-01 00:00:00_diff"
$ws.Range("A15").Value2 = "2023-04-01 00:00:00_diff"
$ws.Range("A16").Value2 = "2023-07-01 00:00:00_diff"
$ws.Range("A17").Value2 = "2023-10-01 00:00:00_diff"
$ws.Range("A18").Value2 = "2024-01-01 00:00:00_diff"
$ws.Range("A19").Value2 = "2024-04-01 00:00:00_diff"
$ws.Range("A20").Value2 = "2024-07-01 00:00:00_diff"
$ws.Range("A21").Value2 = "2024-10-01 00:00:00_diff"
$ws.Range("A22").Value2 = "2025-01-01 00:00:00_diff"
$ws.Range("A23").Value2 = "2025-04-01 00:00:00_diff"
$ws.Range("A24").Value2 = "2025-07-01 00:00:00_diff"

# Data grid B2:K24
# Row 2
$ws.Range("B2").Value2 = 2.249802839716686
$ws.Range("C2").Value2 = 9.718953170716686
$ws.Range("D2").Value2 = -8.452308479283314
$ws.Range("E2").Value2 = -0.3126374792833136
$ws.Range("F2").Value2 = 1.834933388716686
$ws.Range("G2").Value2 = -1.605096479283314
$ws.Range("H2").Value2 = -1.666270479283314
$ws.Range("I2").Value2 = 0.3751855207166864
$ws.Range("J2").Value2 = -0.2074094792833136
$ws.Range("K2").Value2 = -0.1117343352833136
# Row 3
$ws.Range("B3").Value2 = 7.469150331105293
$ws.Range("C3").Value2 = -10.70211131889471
$ws.Range("D3").Value2 = -2.562440318894706
$ws.Range("E3").Value2 = -0.4148694508947059
$ws.Range("F3").Value2 = -3.854899318894706
$ws.Range("G3").Value2 = -3.916073318894706
$ws.Range("H3").Value2 = -1.874617318894706
$ws.Range("I3").Value2 = -2.457212318894706
$ws.Range("J3").Value2 = -2.361537174894706
$ws.Range("K3").Value2 = -2.621956318894706
# Row 4
$ws.Range("B4").Value2 = -18.171261649752
$ws.Range("C4").Value2 = -10.031590649752
$ws.Range("D4").Value2 = -7.884019781751999
$ws.Range("E4").Value2 = -11.324049649752
$ws.Range("F4").Value2 = -11.385223649752
$ws.Range("G4").Value2 = -9.343767649751999
$ws.Range("H4").Value2 = -9.926362649751999
$ws.Range("I4").Value2 = -9.830687505752
$ws.Range("J4").Value2 = -10.091106649752
$ws.Range("K4").Value2 = -9.255651649751998
# Row 5
$ws.Range("B5").Value2 = 8.139671150385469
$ws.Range("C5").Value2 = 10.28724201838547
$ws.Range("D5").Value2 = 6.84721215038547
$ws.Range("E5").Value2 = 6.78603815038547
$ws.Range("F5").Value2 = 8.827494150385469
$ws.Range("G5").Value2 = 8.24489915038547
$ws.Range("H5").Value2 = 8.340574294385469
$ws.Range("I5").Value2 = 8.080155150385469
$ws.Range("J5").Value2 = 8.91561015038547
$ws.Range("K5").Value2 = 8.815537197385471
# Row 6
$ws.Range("B6").Value2 = 2.14757047140632
$ws.Range("C6").Value2 = -1.292459396593679
$ws.Range("D6").Value2 = -1.35363339659368
$ws.Range("E6").Value2 = 0.6878226034063205
$ws.Range("F6").Value2 = 0.1052276034063205
$ws.Range("G6").Value2 = 0.2009027474063205
$ws.Range("H6").Value2 = -0.05951639659367947
$ws.Range("I6").Value2 = 0.7759386034063205
$ws.Range("J6").Value2 = 0.6758656504063205
$ws.Range("K6").Value2 = 0.3220726034063205
# Row 7
$ws.Range("B7").Value2 = -3.440029868393072
$ws.Range("C7").Value2 = -3.501203868393072
$ws.Range("D7").Value2 = -1.459747868393072
$ws.Range("E7").Value2 = -2.042342868393072
$ws.Range("F7").Value2 = -1.946667724393072
$ws.Range("G7").Value2 = -2.207086868393072
$ws.Range("H7").Value2 = -1.371631868393072
$ws.Range("I7").Value2 = -1.471704821393072
$ws.Range("J7").Value2 = -1.825497868393072
$ws.Range("K7").Value2 = -1.677319868393072
# Row 8
$ws.Range("B8").Value2 = -0.06117388186974204
$ws.Range("C8").Value2 = 1.980282118130258
$ws.Range("D8").Value2 = 1.397687118130258
$ws.Range("E8").Value2 = 1.493362262130258
$ws.Range("F8").Value2 = 1.232943118130258
$ws.Range("G8").Value2 = 2.068398118130258
$ws.Range("H8").Value2 = 1.968325165130258
$ws.Range("I8").Value2 = 1.614532118130258
$ws.Range("J8").Value2 = 1.762710118130258
$ws.Range("K8").Value2 = 1.920156118130258
# Row 9
$ws.Range("B9").Value2 = 2.041456296544459
$ws.Range("C9").Value2 = 1.458861296544459
$ws.Range("D9").Value2 = 1.554536440544459
$ws.Range("E9").Value2 = 1.294117296544459
$ws.Range("F9").Value2 = 2.129572296544459
$ws.Range("G9").Value2 = 2.029499343544459
$ws.Range("H9").Value2 = 1.675706296544459
$ws.Range("I9").Value2 = 1.823884296544459
$ws.Range("J9").Value2 = 1.981330296544459
$ws.Range("K9").Value2 = 1.480667296544459
# Row 10
$ws.Range("B10").Value2 = -0.5825946406117954
$ws.Range("C10").Value2 = -0.4869194966117955
$ws.Range("D10").Value2 = -0.7473386406117954
$ws.Range("E10").Value2 = 0.08811635938820461
$ws.Range("F10").Value2 = -0.01195659361179541
$ws.Range("G10").Value2 = -0.3657496406117954
$ws.Range("H10").Value2 = -0.2175716406117954
$ws.Range("I10").Value2 = -0.06012564061179543
$ws.Range("J10").Value2 = -0.5607886406117955
$ws.Range("K10").Value2 = -0.2804276406117954
# Row 11
$ws.Range("B11").Value2 = 0.09567504042184558
$ws.Range("C11").Value2 = -0.1647441035781544
$ws.Range("D11").Value2 = 0.6707108964218456
$ws.Range("E11").Value2 = 0.5706379434218456
$ws.Range("F11").Value2 = 0.2168448964218456
$ws.Range("G11").Value2 = 0.3650228964218456
$ws.Range("H11").Value2 = 0.5224688964218456
$ws.Range("I11").Value2 = 0.02180589642184558
$ws.Range("J11").Value2 = 0.3021668964218456
$ws.Range("K11").Value2 = 0.1301918964218456
# Row 12
$ws.Range("B12").Value2 = -0.2604191443875122
$ws.Range("C12").Value2 = 0.5750358556124878
$ws.Range("D12").Value2 = 0.4749629026124878
$ws.Range("E12").Value2 = 0.1211698556124878
$ws.Range("F12").Value2 = 0.2693478556124878
$ws.Range("G12").Value2 = 0.4267938556124878
$ws.Range("H12").Value2 = -0.07386914438751221
$ws.Range("I12").Value2 = 0.2064918556124878
$ws.Range("J12").Value2 = 0.0345168556124878
$ws.Range("K12").Value2 = 0.3403798556124878
# Row 13
$ws.Range("B13").Value2 = 0.8354548926112106
$ws.Range("C13").Value2 = 0.7353819396112107
$ws.Range("D13").Value2 = 0.3815888926112106
$ws.Range("E13").Value2 = 0.5297668926112107
$ws.Range("F13").Value2 = 0.6872128926112107
$ws.Range("G13").Value2 = 0.1865498926112106
$ws.Range("H13").Value2 = 0.4669108926112106
$ws.Range("I13").Value2 = 0.2949358926112106
$ws.Range("J13").Value2 = 0.6007988926112107
$ws.Range("K13").Value2 = -0.01153110738878937
# Row 14
$ws.Range("B14").Value2 = -0.1000730565472806
$ws.Range("C14").Value2 = -0.4538661035472806
$ws.Range("D14").Value2 = -0.3056881035472806
$ws.Range("E14").Value2 = -0.1482421035472806
$ws.Range("F14").Value2 = -0.6489051035472806
$ws.Range("G14").Value2 = -0.3685441035472806
$ws.Range("H14").Value2 = -0.5405191035472806
$ws.Range("I14").Value2 = -0.2346561035472806
$ws.Range("J14").Value2 = -0.8469861035472805
$ws.Range("K14").Value2 = -0.1586151035472806
# Row 15
$ws.Range("B15").Value2 = -0.3537867436446591
$ws.Range("C15").Value2 = -0.2056087436446591
$ws.Range("D15").Value2 = -0.04816274364465911
$ws.Range("E15").Value2 = -0.5488257436446591
$ws.Range("F15").Value2 = -0.2684647436446591
$ws.Range("G15").Value2 = -0.4404397436446591
$ws.Range("H15").Value2 = -0.1345767436446591
$ws.Range("I15").Value2 = -0.746906743644659
$ws.Range("J15").Value2 = -0.05853574364465908
# Row 16
$ws.Range("B16").Value2 = 0.1481777624350372
$ws.Range("C16").Value2 = 0.3056237624350371
$ws.Range("D16").Value2 = -0.1950392375649629
$ws.Range("E16").Value2 = 0.08532176243503714
$ws.Range("F16").Value2 = -0.08665323756496286
$ws.Range("G16").Value2 = 0.2192097624350371
$ws.Range("H16").Value2 = -0.3931202375649628
$ws.Range("I16").Value2 = 0.2952507624350372
# Row 17
$ws.Range("B17").Value2 = 0.1574463720025918
$ws.Range("C17").Value2 = -0.3432166279974082
$ws.Range("D17").Value2 = -0.06285562799740819
$ws.Range("E17").Value2 = -0.2348306279974082
$ws.Range("F17").Value2 = 0.0710323720025918
$ws.Range("G17").Value2 = -0.5412976279974082
$ws.Range("H17").Value2 = 0.1470733720025918
# Row 18
$ws.Range("B18").Value2 = -0.5006626170015632
$ws.Range("C18").Value2 = -0.2203016170015632
$ws.Range("D18").Value2 = -0.3922766170015632
$ws.Range("E18").Value2 = -0.08641361700156319
$ws.Range("F18").Value2 = -0.6987436170015632
$ws.Range("G18").Value2 = -0.01037261700156317
# Row 19
$ws.Range("B19").Value2 = 0.2803578395245076
$ws.Range("C19").Value2 = 0.1083828395245076
$ws.Range("D19").Value2 = 0.4142458395245076
$ws.Range("E19").Value2 = -0.1980841604754924
$ws.Range("F19").Value2 = 0.4902868395245076
# Row 20
$ws.Range("B20").Value2 = -0.1719750410109616
$ws.Range("C20").Value2 = 0.1338879589890384
$ws.Range("D20").Value2 = -0.4784420410109615
$ws.Range("E20").Value2 = 0.2099289589890385
# Row 21
$ws.Range("B21").Value2 = 0.3058628168340501
$ws.Range("C21").Value2 = -0.3064671831659499
$ws.Range("D21").Value2 = 0.3819038168340501
# Row 22
$ws.Range("B22").Value2 = -0.6123297229122814
$ws.Range("C22").Value2 = 0.07604127708771863
# Row 23
$ws.Range("B23").Value2 = 0.6883712297750049

Write-Host "edit complete"